$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated Price (D) and Volume(1h) (E) values scraped for this run.
# D-column cells are forced to Text format first so purely-numeric-looking
# price strings (e.g. "243.52") are not auto-converted to Excel numbers,
# then the style is reset to Normal so no stray formatting is introduced.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.578.37"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.10%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.961.94"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.70%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.628"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.33"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.73%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  +3.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0786"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.29%  "
$ws.Range("E11").Value = "  +0.56%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.17"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.842"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.247.94"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.56"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.30"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.37%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.963.59"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.55%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "36.510.03"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.24"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0853"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.49%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "229.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.49%  "
$ws.Range("E22").Value = "  +1.08%  "
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("E24").Value = "  +1.31%  "
$ws.Range("E25").Value = "  +1.98%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.144"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +7.77%  "
$ws.Range("E27").Value = "  -1.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "161.14"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.27"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.37%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.33"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +20.66%  "
$ws.Range("E31").Value = "  +1.69%  "
$ws.Range("E32").Value = "  +2.47%  "
$ws.Range("E33").Value = "  -1.25%  "
$ws.Range("E34").Value = "  +6.50%  "
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("E36").Value = "  +2.60%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.38"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.82%  "
$ws.Range("E38").Value = "  -0.06%  "
$ws.Range("E39").Value = "  -13.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0970"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.82%  "
$ws.Range("E41").Value = "  +0.23%  "
$ws.Range("E42").Value = "  +0.62%  "
$ws.Range("E43").Value = "  -0.99%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.88"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.364.57"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.59%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "88.67"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.94%  "
$ws.Range("E47").Value = "  -0.73%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.23"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.95%  "
$ws.Range("E49").Value = "  +0.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "45.93"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.144.11"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.77%  "
